# Journal de travail - add new time entries for the week of 2017-04-03 (rows 36-42)
# and remove now-unneeded blank rows so totals line up again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# ------------------------------------------------------------------
# 1) Delete 5 now-superfluous blank rows (were rows 43-47) so the
#    "Total" row of this week's block moves back up against the data.
# ------------------------------------------------------------------
$ws.Range("A43:A47").EntireRow.Delete()

# ------------------------------------------------------------------
# 2) Fill in the week's time entries (rows 36-42)
# ------------------------------------------------------------------
$ws.Range("C36").Value = 0.33333333333333331
$ws.Range("D36").Value = "-"
$ws.Range("E36").Value = 0.41666666666666669
$ws.Range("F36").Value = "Validation De la vue createShoot"

$ws.Range("C37").Value = 0.41666666666666669
$ws.Range("D37").Value = "-"
$ws.Range("E37").Value = 0.5
$ws.Range("F37").Value = "Création du composant arrowItem et intégration dans le editShoot"

$ws.Range("C38").Value = 0.53125
$ws.Range("D38").Value = "-"
$ws.Range("E38").Value = 0.5625
$ws.Range("F38").Value = "Reflexions sur la praticité des ajouts de flèches"

$ws.Range("C39").Value = 0.5625
$ws.Range("D39").Value = "-"
$ws.Range("E39").Value = 0.625
$ws.Range("F39").Value = "Mise en place de l'ajout des flèches dans la Vue EditShoot"

$ws.Range("C40").Value = 0.625
$ws.Range("D40").Value = "-"
$ws.Range("E40").Value = 0.66666666666666663
$ws.Range("F40").Value = "Ajout de paramétres dans la route editShoot pour pas avoir écran blanc à l'actualisation"
$ws.Rows.Item(40).RowHeight = 32.25

$ws.Range("C41").Value = 0.66666666666666663
$ws.Range("D41").Value = "-"
$ws.Range("E41").Value = 0.6875
$ws.Range("F41").Value = "Aide du chef de projet afin d'ordonner les Arrows lors de l'editing"

$ws.Range("C42").Value = 0.6875
$ws.Range("D42").Value = "-"
$ws.Range("E42").Value = 0.71458333333333324
$ws.Range("F42").Value = "Rédaction rapport partie editShoot"

# ------------------------------------------------------------------
# 3) Restore the week's "Total" formula in row 43 (lost its terms
#    when the blank rows above were removed)
# ------------------------------------------------------------------
$ws.Range("E43").Formula = "=E42-C42+E41-C41+E40-C40+E39-C39+E38-C38+E37-C37+E36-C36"

# ------------------------------------------------------------------
# 4) Make sure the grand-total formula (H5) references the
#    up-to-date block totals
# ------------------------------------------------------------------
$ws.Range("H5").Formula = "=SUM(E8,E14,E22,E29,E43,E50,E59,E76,E80,E35,E65,E68)*24"

# ------------------------------------------------------------------
# 5) Shrink the print area now that the sheet is 5 rows shorter
# ------------------------------------------------------------------
$wb.Names.Item(1).RefersTo = "='Journal de travail'!`$A`$1:`$I`$80"

# ------------------------------------------------------------------
# 6) Refresh the view state to match what was left active
# ------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("C44:E44").Select()

$wb.Application.CalculateFull()
